# Service - NodePort - cmds (2)
#
# Each slide (slide 1 and slide 2) has a "Service" YAML code textbox that
# currently lists:
#     - port: 80
#       protocol: TCP
#       targetPort: 80
#
# It needs to become:
#     - port: 80
#       targetPort: 80
#       nodePort: 30008
#
# i.e. the "protocol: TCP" line is turned into "targetPort: 80" (key text +
# value text + value color change, since the value becomes a number), and
# the old "targetPort: 80" line is turned into "nodePort: 30008" (key text +
# value text only, same colors as before).

$p = $ppt.ActivePresentation

function Update-ServiceYamlShape {
    param($shape)

    $tr = $shape.TextFrame.TextRange

    # Paragraph 9: "    protocol: TCP"  ->  "    targetPort: 80"
    $para = $tr.Paragraphs(9, 1)

    $keyRun = $para.Runs(2, 1)
    $keyRun.Text = "targetPort"

    $valueRun = $para.Runs(4, 1)
    $valueRun.Text = "80"
    $valueRun.Font.Color.RGB = 0xA8CEB5   # COM BGR for srgbClr B5CEA8

    # Paragraph 10: "    targetPort: 80"  ->  "    nodePort: 30008"
    $para2 = $tr.Paragraphs(10, 1)

    $keyRun2 = $para2.Runs(2, 1)
    $keyRun2.Text = "nodePort"

    $valueRun2 = $para2.Runs(4, 1)
    $valueRun2.Text = "30008"
}

# Slide 1 -> shape "Title 1" (index 1) holds the Service yaml block.
$slide1 = $p.Slides.Item(1)
Update-ServiceYamlShape $slide1.Shapes.Item(1)

# Slide 2 -> shape "Title 1" (index 2) holds the Service yaml block.
$slide2 = $p.Slides.Item(2)
Update-ServiceYamlShape $slide2.Shapes.Item(2)
